$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 400
    5  = 450
    6  = 800
    7  = 800
    9  = 700
    10 = 666.67
    11 = 300
    12 = 300
    13 = 800
    14 = 800
    15 = 800
    18 = 300
    20 = 400
    21 = 300
    22 = 800
    23 = 800
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
